# 4-06-2024 work done in half day
# Adds the task / status / remark for 4-06-2024 (row 30) and scrolls /
# selects the relevant cell, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task text, status ("no") and remark ("half day") for the 4-06-2024 entry.
$ws.Range("C30").Value = "make profile menu list profile and logout and add icon on dashboard and fixed ipad mini sized "
$ws.Range("D30").Value = "no"
$ws.Range("E30").Value = "half day"

# Reflect the author's final cursor position / scroll state in the sheet view.
$ws.Range("E30").Select()

$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$win.ScrollRow = 18
